$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab
$ws.Name = "Control de Gastos"

# Row 1
$ws.Range("A1").Value = "Saldo inicial:"
$ws.Range("B1").Value = 1000

# Row 2
$ws.Range("A2").Value = " "

# Row 3 - headers
$ws.Range("A3").Value = "Ingresos"
$ws.Range("B3").Value = "Gastos"
$ws.Range("C3").Value = "Fecha"

# Row 4
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = 60
$ws.Range("C4").Value = 44087
$ws.Range("C4").NumberFormat = "yyyy-mm-dd"
